# Adding files to 29th Sep folder
#
# On slide 4, the "Picture 2" picture (the one sitting above the
# "Straight Connector 17" line, at x=1666070 EMU) is moved up, and a
# copy of it is dropped in further down the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# EMU -> Points helper. PowerPoint's Left/Top/Width/Height are in points
# (1 pt = 12700 EMU). The COM bridge stores these as single-precision
# floats, so a plain division can truncate one EMU short once the value
# is converted back on save (e.g. 4336837 -> 4336836). Nudge by a tiny
# fraction of a point (way below a hundredth of an EMU's worth of
# visual difference) so the round-trip lands back on the exact EMU.
function EMUToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00002
}

# Locate the picture shape we need to move: the "Picture 2" shape
# (id=2050) positioned at x=1666070, y=3051292 EMU.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 2050) {
        $target = $sh
        break
    }
}
if ($target -eq $null) {
    # Fallback: match on name + approximate position instead of Id.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -eq "Picture 2" -and [Math]::Round($sh.Left * 12700) -eq 1666070) {
            $target = $sh
            break
        }
    }
}

# Duplicate it in place before moving it, so the new copy keeps all of
# the original picture's formatting (picture locks, hidden-fill ext,
# srcRect, etc.) and its image relationship.
$newPic = $target.Duplicate()

# Move the new copy to its final resting spot further down the slide.
$newPic.Left = EMUToPt 1666070
$newPic.Top = EMUToPt 4336837

# Move the original picture up.
$target.Left = EMUToPt 1666070
$target.Top = EMUToPt 1118797
